# Clean up form formatting
# - Make the "survey" sheet the active/selected sheet (was "settings").
# - Update the selected cell in the survey sheet's frozen bottom-right pane.
# - Convert the two "hidden" question rows to "text" questions with an
#   explicit "hidden" appearance column, and the "string" question type to
#   "text" (XLSForm style cleanup).

$wb = $excel.ActiveWorkbook

$survey = $wb.Worksheets.Item("survey")

# --- Row 4: source (was type=hidden) -> type=text, appearance=hidden -----
$survey.Range("A4").Value = "text"
$survey.Range("F4").Value = "hidden"

# --- Row 5: source_id (was type=hidden) -> type=text, appearance=hidden --
$survey.Range("A5").Value = "text"
$survey.Range("F5").Value = "hidden"

# --- Row 7: _id (was type=string) -> type=text ----------------------------
$survey.Range("A7").Value = "text"

# --- Make "survey" the active sheet and restore its selection ------------
$survey.Activate() | Out-Null
$survey.Range("C19").Select() | Out-Null

# "settings" is no longer the active sheet; keep its own saved selection
# (E8) untouched.
